# "Keep it simple, keep it stupid." + template for new vs 2008 project
#
# Adds a "Serial" (column B) family label to a number of H.26x / MPEG rows
# on Sheet1, fixes a typo in the JPEG-standard description, and leaves the
# selection on E14 (matching the author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# H.26x family rows
$ws.Range("B3").Value  = "H.26x"   # 1984 H.120
$ws.Range("B4").Value  = "H.26x"   # 1988 H.261
$ws.Range("B13").Value = "H.26x"   # 1995 H.262
$ws.Range("B31").Value = "H.26x"   # 2003 H.264
$ws.Range("B34").Value = "H.26x"   # 2004 H.265
$ws.Range("B44").Value = "H.26x"   # 2010 H.265

# MPEG family rows
$ws.Range("B5").Value  = "MPEG"    # 1988 MPEG-1
$ws.Range("B7").Value  = "MPEG"    # 1990 MPEG-2
$ws.Range("B11").Value = "MPEG"    # 1993 MPEG-1
$ws.Range("B12").Value = "MPEG"    # 1994 MPEG-2
$ws.Range("B14").Value = "MPEG"    # 1996 MPEG-2
$ws.Range("B19").Value = "MPEG"    # 1998 MPEG-2
$ws.Range("B20").Value = "MPEG"    # 1998 MPEG-4
$ws.Range("B21").Value = "MPEG"    # 1999 MPEG-4
$ws.Range("B23").Value = "MPEG"    # 2000 MPEG-2
$ws.Range("B26").Value = "MPEG"    # 2001 DivX Encore 2
$ws.Range("B29").Value = "MPEG"    # 2003 MPEG-4
$ws.Range("B39").Value = "MPEG"    # 2007 MPEG-H
$ws.Range("B42").Value = "MPEG"    # 2009 MPEG-H
$ws.Range("B48").Value = "MPEG"    # 2013 MPEG-2

# Fix typo "puublished" -> "published"
$ws.Range("E10").Value = "Firsh published JPEG standard"

# Leave the selection where the author left it
$ws.Range("E14").Select()
